# Master user belum fix
# Add a new data row (row 10) to Sheet1 with the "Gedung Lab Fisika" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Gedung Lab Fisika"
$ws.Range("C10").Value = 67.8
$ws.Range("D10").Value = 89.8
$ws.Range("E10").Value = "Legalitas"
$ws.Range("F10").Value = "Negara"
$ws.Range("G10").Value = "Hak A"
$ws.Range("H10").Value = 1714.98
$ws.Range("I10").Value = 7
$ws.Range("J10").Value = 9000.0
$ws.Range("K10").Value = 899.0
$ws.Range("L10").Value = "Sedang"
$ws.Range("M10").Value = "Khusus"
$ws.Range("N10").Value = "Lokasi Kepadatan Sedang"
$ws.Range("O10").Value = "Permanen"
$ws.Range("P10").Value = "Rendah"
$ws.Range("Q10").Value = "Pasif"
$ws.Range("R10").Value = "i"
$ws.Range("S10").Value = "j"
$ws.Range("T10").Value = "k"
